# "Generate Report for Handback" - refresh the timestamps recorded the last
# time the handback status report was generated.
#
#   Overview!G2 (Latest HO Xliff Generate Date)         -> 2016-08-16 00:59:29
#   zh-cn!H2   (Correspond Handoff Datetime, row 2)      -> 2016-08-16 00:59:24
#   zh-cn!K2   (Correspond Handback DateTime, row 2)     -> 2016-08-16 00:59:41
#   de-de!H2   (Correspond Handoff Datetime, row 2)      -> 2016-08-16 00:59:29
#   de-de!K2   (Correspond Handback DateTime, row 2)     -> 2016-08-16 00:59:48

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-16 00:59:29"

$wsZhCn.Range("H2").Value = "2016-08-16 00:59:24"
$wsZhCn.Range("K2").Value = "2016-08-16 00:59:41"

$wsDeDe.Range("H2").Value = "2016-08-16 00:59:29"
$wsDeDe.Range("K2").Value = "2016-08-16 00:59:48"
